# Add a new "2021" column (R) to the table, mirroring the formatting used
# by the existing 2020 column (Q), then move the active selection to T3 -
# matching the saved workbook view recorded after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing column (Q) into the new column
# (R) for every row of the table, then fill in the new values.

# Row 2 - thin bottom-bordered spacer row above the header (no value)
$ws.Range("Q2").Copy()
$ws.Range("R2").PasteSpecial(-4122)

# Row 3 - year header row
$ws.Range("Q3").Copy()
$ws.Range("R3").PasteSpecial(-4122)
$ws.Range("R3").Value = 2021

# Row 4 - GVA share of manufacturing output in GDP, %
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)
$ws.Range("R4").Value = 13.5

# Row 5 - GVA of manufacturing industry in GDP per capita
$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122)
$ws.Range("R5").Value = 15.1

$excel.CutCopyMode = $false

# Move/update the active selection as recorded in the saved workbook view.
$ws.Range("T3").Select()
